$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corregir errores de Precedencias
# C9: precedencia combinada "4, 5" (antes 4.5 como numero)
$ws.Range("C9").Value = "4, 5"

# C10: precedencia corregida de 7 a 6
$ws.Range("C10").Value = 6

# D11: duracion corregida de 1 a 3
$ws.Range("D11").Value = 3

# Ajustar la vista: desplazar a la celda superior izquierda B1 y seleccionar D11
$ws.Range("D11").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
